# Add season-record columns (Wins/Losses/Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - style matches the other header cells (bold, bordered, centered)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-40: constant record for every player row
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 72   # AD
    $ws.Cells.Item($r, 31).Value = 89   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
